$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user_h")

# Add new row 33 with Mac-Address / Document Type entry (mirrors existing rows' pattern)
$ws.Cells.Item(33, 1).Value = 10002
$ws.Cells.Item(33, 2).Value = 110032
$ws.Cells.Item(33, 3).Value = "eng"
$ws.Cells.Item(33, 4).Value = $true
$ws.Cells.Item(33, 5).Value = "superadmin"
$ws.Cells.Item(33, 6).Value = "now()"
$ws.Cells.Item(33, 7).Value = "now()"

# Update the selected cell to match the author's final cursor position
$ws.Range("D26").Select()
